$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2538.8809
$ws.Range("I15").Value = 2538.8809
$ws.Range("K15").Value = 7616.6427
$ws.Range("M15").Value = -7447.6427

$ws.Range("H70").Value = 1063.2632
$ws.Range("I70").Value = 883.55554
$ws.Range("J70").Value = 1225
$ws.Range("K70").Value = 2650.66662
$ws.Range("L70").Value = 3675
$ws.Range("M70").Value = -2380.66662
$ws.Range("N70").Value = -4215

$ws.Range("H73").Value = 1063.2632
$ws.Range("I73").Value = 883.55554
$ws.Range("J73").Value = 1225
$ws.Range("K73").Value = 2650.66662
$ws.Range("L73").Value = 3675
$ws.Range("M73").Value = -1714.66662
$ws.Range("N73").Value = -5547

$ws.Range("H75").Value = 42000
$ws.Range("J75").Value = 42000
$ws.Range("L75").Value = 42000
$ws.Range("N75").Value = -43872

$ws.Range("H78").Value = 42000
$ws.Range("J78").Value = 42000
$ws.Range("L78").Value = 126000
$ws.Range("N78").Value = -135360

$ws.Range("H107").Value = 50654.75
$ws.Range("J107").Value = 922.5
$ws.Range("L107").Value = 922.5
$ws.Range("N107").Value = -4762.5

$ws.Range("H113").Value = 29415524
$ws.Range("I113").Value = 71430130
$ws.Range("J113").Value = 5298.7
$ws.Range("K113").Value = 71430130
$ws.Range("L113").Value = 5298.7
$ws.Range("M113").Value = -71426876
$ws.Range("N113").Value = -11806.7

$ws.Range("H116").Value = 3765.4375
$ws.Range("I116").Value = 2150
$ws.Range("J116").Value = 6457.8335
$ws.Range("K116").Value = 2150
$ws.Range("L116").Value = 6457.8335
$ws.Range("M116").Value = 1292
$ws.Range("N116").Value = -13341.8335

$ws.Range("H124").Value = 48320
$ws.Range("J124").Value = 48320
$ws.Range("L124").Value = 48320
$ws.Range("N124").Value = -58140

$ws.Range("H130").Value = 51324.445
$ws.Range("J130").Value = 51324.445
$ws.Range("L130").Value = 51324.445
$ws.Range("N130").Value = -61364.445

$ws.Range("H132").Value = 3460.25
$ws.Range("I132").Value = 2319.2917
$ws.Range("J132").Value = 5742.1665
$ws.Range("K132").Value = 6957.875100000001
$ws.Range("L132").Value = 17226.4995
$ws.Range("M132").Value = -4427.875100000001
$ws.Range("N132").Value = -22286.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23103.834
$ws.Range("I32").Value = 4286.352
$ws.Range("J32").Value = 107782.5
$ws.Range("K32").Value = 4286.352
$ws.Range("L32").Value = 107782.5
$ws.Range("M32").Value = -3999.352
$ws.Range("N32").Value = -108356.5

$ws.Range("H80").Value = 26577.111
$ws.Range("J80").Value = 26577.111
$ws.Range("L80").Value = 26577.111
$ws.Range("N80").Value = -28573.111

$ws.Range("H83").Value = 26577.111
$ws.Range("J83").Value = 26577.111
$ws.Range("L83").Value = 79731.333
$ws.Range("N83").Value = -89715.333

$ws.Range("H125").Value = 112543230
$ws.Range("J125").Value = 112543230
$ws.Range("L125").Value = 112543230
$ws.Range("N125").Value = -112553070

$ws.Range("H129").Value = 42599.6
$ws.Range("J129").Value = 42599.6
$ws.Range("L129").Value = 42599.6
$ws.Range("N129").Value = -52599.6

$ws.Range("H131").Value = 44979.5
$ws.Range("J131").Value = 44979.5
$ws.Range("L131").Value = 44979.5
$ws.Range("N131").Value = -55059.5

$ws.Range("H132").Value = 121103.734
$ws.Range("I132").Value = 147963.5
$ws.Range("J132").Value = 6949.75
$ws.Range("K132").Value = 443890.5
$ws.Range("L132").Value = 20849.25
$ws.Range("M132").Value = -441360.5
$ws.Range("N132").Value = -25909.25

$ws.Range("H133").Value = 50261
$ws.Range("J133").Value = 50261
$ws.Range("L133").Value = 50261
$ws.Range("N133").Value = -55321

$ws.Range("H134").Value = 47872.168
$ws.Range("J134").Value = 47872.168
$ws.Range("L134").Value = 47872.168
$ws.Range("N134").Value = -58012.168

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18512.715
$ws.Range("J82").Value = 19888.666
$ws.Range("L82").Value = 19888.666
$ws.Range("N82").Value = -20654.666

$ws.Range("H85").Value = 18512.715
$ws.Range("J85").Value = 19888.666
$ws.Range("L85").Value = 19888.666
$ws.Range("N85").Value = -22540.666

$ws.Range("H122").Value = 49882.855
$ws.Range("J122").Value = 49882.855
$ws.Range("L122").Value = 49882.855
$ws.Range("N122").Value = -59682.855

$ws.Range("H132").Value = 39000.8
$ws.Range("J132").Value = 39000.8
$ws.Range("L132").Value = 39000.8
$ws.Range("N132").Value = -49120.8

$ws.Range("H135").Value = 53780
$ws.Range("J135").Value = 53780
$ws.Range("L135").Value = 53780
$ws.Range("N135").Value = -63920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 20654
$ws.Range("J109").Value = 20654
$ws.Range("L109").Value = 20654
$ws.Range("N109").Value = -22734

$ws.Range("H123").Value = 52653.332
$ws.Range("J123").Value = 52653.332
$ws.Range("L123").Value = 52653.332
$ws.Range("N123").Value = -62453.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 7363.6113
$ws.Range("I64").Value = 812
$ws.Range("J64").Value = 7749
$ws.Range("K64").Value = 2436
$ws.Range("L64").Value = 23247
$ws.Range("M64").Value = -2166
$ws.Range("N64").Value = -23787

$ws.Range("H67").Value = 7363.6113
$ws.Range("I67").Value = 812
$ws.Range("J67").Value = 7749
$ws.Range("K67").Value = 2436
$ws.Range("L67").Value = 23247
$ws.Range("M67").Value = -1500
$ws.Range("N67").Value = -25119

$ws.Range("H87").Value = 8864.909
$ws.Range("I87").Value = 5502
$ws.Range("J87").Value = 14750
$ws.Range("K87").Value = 16506
$ws.Range("L87").Value = 44250
$ws.Range("M87").Value = -15258
$ws.Range("N87").Value = -46746

$ws.Range("H90").Value = 8864.909
$ws.Range("I90").Value = 5502
$ws.Range("J90").Value = 14750
$ws.Range("K90").Value = 49518
$ws.Range("L90").Value = 132750
$ws.Range("M90").Value = -43278
$ws.Range("N90").Value = -145230

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 6150.8335
$ws.Range("J43").Value = 9435
$ws.Range("L43").Value = 9435
$ws.Range("N43").Value = -9737

$ws.Range("H57").Value = 17072.7
$ws.Range("J57").Value = 17072.7
$ws.Range("L57").Value = 17072.7
$ws.Range("N57").Value = -18712.7

$ws.Range("H124").Value = 49980
$ws.Range("J124").Value = 49980
$ws.Range("L124").Value = 49980
$ws.Range("N124").Value = -59800

$ws.Range("H127").Value = 46036
$ws.Range("J127").Value = 46036
$ws.Range("L127").Value = 46036
$ws.Range("N127").Value = -55956

$ws.Range("H128").Value = 51932.5
$ws.Range("J128").Value = 51932.5
$ws.Range("L128").Value = 51932.5
$ws.Range("N128").Value = -61892.5

$ws.Range("H134").Value = 14900
$ws.Range("J134").Value = 14900
$ws.Range("L134").Value = 44700
$ws.Range("N134").Value = -49770

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 5600
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 5600
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 5600
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -5824

$ws.Range("H30").Value = 18855.334
$ws.Range("I30").Value = 18626.4
$ws.Range("J30").Value = 20000
$ws.Range("K30").Value = 18626.4
$ws.Range("L30").Value = 20000
$ws.Range("M30").Value = -18518.4
$ws.Range("N30").Value = -20216

$ws.Range("H109").Value = 28250
$ws.Range("J109").Value = 28250
$ws.Range("L109").Value = 28250
$ws.Range("N109").Value = -31024

$ws.Range("H127").Value = 42050.715
$ws.Range("J127").Value = 42050.715
$ws.Range("L127").Value = 42050.715
$ws.Range("N127").Value = -51970.715

$ws.Range("H128").Value = 52108
$ws.Range("J128").Value = 52108
$ws.Range("L128").Value = 52108
$ws.Range("N128").Value = -62068

$ws.Range("H129").Value = 48482.25
$ws.Range("J129").Value = 48482.25
$ws.Range("L129").Value = 48482.25
$ws.Range("N129").Value = -58482.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7500
$ws.Range("J41").Value = 7500
$ws.Range("L41").Value = 7500
$ws.Range("N41").Value = -8280

$ws.Range("H125").Value = 49361.25
$ws.Range("J125").Value = 49361.25
$ws.Range("L125").Value = 49361.25
$ws.Range("N125").Value = -59201.25

$ws.Range("H129").Value = 26254
$ws.Range("J129").Value = 26254
$ws.Range("L129").Value = 26254
$ws.Range("N129").Value = -36254
